$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT (e.g. numeric-looking price strings)
# so Excel does not auto-convert it to a Double and lose formatting/precision.
function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue "D2" "244.71"
Set-TextValue "D3" "23.15"
Set-TextValue "D4" "5.438"
Set-TextValue "D6" "3.389"
Set-TextValue "D7" "0.8092"
Set-TextValue "D8" "0.9259"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D9" "0.01119"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1428"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07423"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03384"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03038"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09349"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.952"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001602"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04824"
$ws.Range("E17").Value = "16CoinExTokenCET"
Set-TextValue "D18" "0.005535"
Set-TextValue "D19" "0.004153"
Set-TextValue "D20" "0.0009840"
Set-TextValue "D21" "0.00007706"
Set-TextValue "D22" "3.660"
Set-TextValue "D23" "6.454"
Set-TextValue "D26" "0.1304"
Set-TextValue "D41" "0.006215"
Set-TextValue "D43" "0.002902"
Set-TextValue "D44" "0.007235"
$ws.Range("E44").Value = "43LocalTradersLCT"
Set-TextValue "D45" "0.00005144"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.0005803"
Set-TextValue "D48" "0.8556"
Set-TextValue "D50" "0.00002102"
Set-TextValue "D51" "0.0002002"
